# Generate Report for Archive
#
# The localization status for the reported file moved on from "Ready for
# handoff" to "In Translation". That status string is shown on the
# "Overview" sheet (once per target language, columns E and F) and on each
# per-language sheet ("zh-cn", "de-de") in the "Status" column (C).
#
# Updating the text also changes the natural (content-based) width of
# those status columns, so we shrink them to match the new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-language sheets: "Status" column (C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
